# Re-apply the R script's recalculated adjusted-close series (a dividend was
# incorporated, so the full back-adjustment history for IMP.MI shifts by a
# uniform factor) and append the newly scraped trading day as row 482.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# adj_close (column G) values for every data row, in order: existing rows
# 2..481 first, then the brand-new row 482 appended at the end.
$adjClose = @("4.51669645309448","4.50708675384521","4.4590368270874","4.4590368270874","4.4590368270874","4.57435655593872","4.57435655593872","4.70889663696289","4.65123701095581","4.65123701095581","4.65123701095581","4.65123701095581","4.70889663696289","4.78097152709961","4.6896767616272","4.6896767616272","4.6896767616272","4.70889663696289","4.70889663696289","4.70889663696289","4.70889663696289","4.70889663696289","4.70889663696289","4.70889663696289","4.70889663696289","4.80499649047852","4.80499649047852","4.61279678344727","4.43981695175171","4.43981695175171","4.43981695175171","4.43981695175171","4.43981695175171","4.43981695175171","4.43981695175171","4.43981695175171","4.46864700317383","4.30527687072754","4.30527687072754","4.30527687072754","4.30527687072754","4.30527687072754","4.30527687072754","4.30527687072754","4.30527687072754","4.30527687072754","4.17554187774658","4.17554187774658","4.17554187774658","4.17554187774658","4.17554187774658","4.17554187774658","4.17554187774658","4.17554187774658","4.17554187774658","4.03619718551636","4.03619718551636","4.03619718551636","4.03619718551636","3.63257741928101","3.63257741928101","3.94009709358215","4.03619718551636","3.97853708267212","4.13229703903198","4.20917701721191","4.28605651855469","4.28605651855469","4.28605651855469","4.28605651855469","4.28605651855469","4.28605651855469","4.40137672424316","4.40137672424316","4.42059659957886","4.42059659957886","4.53591632843018","4.53591632843018","4.4782567024231","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.63201665878296","4.4782567024231","4.4782567024231","4.42059659957886","4.42059659957886","4.36293649673462","4.36293649673462","4.36293649673462","4.36293649673462","4.36293649673462","4.36293649673462","4.36293649673462","4.36293649673462","4.4590368270874","4.4590368270874","4.38215684890747","4.38215684890747","4.38215684890747","4.38215684890747","4.38215684890747","4.38215684890747","4.38215684890747","4.38215684890747","4.38215684890747","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.24761724472046","4.34371662139893","4.24761724472046","4.24761724472046","4.24761724472046","4.32449674606323","4.32449674606323","4.32449674606323","4.32449674606323","4.32449674606323","4.32449674606323","4.32449674606323","4.32449674606323","4.42059659957886","4.55513620376587","4.55513620376587","4.55513620376587","4.90109634399414","4.63201665878296","4.55513620376587","4.22839689254761","4.15151691436768","4.07463645935059","4.01697683334351","4.01697683334351","4.01697683334351","3.97853708267212","3.97853708267212","3.97853708267212","3.97853708267212","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.82477712631226","3.7478973865509","3.7478973865509","3.7478973865509","3.7478973865509","3.7478973865509","3.7478973865509","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.86321711540222","3.76711750030518","3.76711750030518","3.76711750030518","3.6517972946167","3.6517972946167","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.55569744110107","3.45959758758545","3.45959758758545","3.51725745201111","3.51725745201111","3.51725745201111","3.5801465511322","3.59960389137268","3.59960389137268","3.59960389137268","3.59960389137268","3.48285984992981","3.5801465511322","3.52177429199219","3.52177429199219","3.52177429199219","3.52177429199219","3.52177429199219","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.5801465511322","3.46340250968933","3.46340250968933","3.46340250968933","3.46340250968933","3.46340250968933","3.46340250968933","3.4050304889679","3.4050304889679","3.4050304889679","3.4050304889679","3.4050304889679","3.4050304889679","3.4050304889679","3.4050304889679","3.36611604690552","3.59960389137268","3.59960389137268","3.59960389137268","3.59960389137268","3.59960389137268","3.63851833343506","3.65797567367554","3.65797567367554","3.65797567367554","3.65797567367554","3.65797567367554","3.63851833343506","3.63851833343506","3.63851833343506","3.59960389137268","3.59960389137268","3.59960389137268","3.59960389137268","3.59960389137268","4.06657934188843","3.69689035415649","3.73580479621887","3.8914635181427","3.8914635181427","3.8914635181427","3.8914635181427","3.8914635181427","3.8914635181427","3.8914635181427","3.73580479621887","3.73580479621887","3.73580479621887","3.63851833343506","3.69689035415649","3.69689035415649","3.69689035415649","3.69689035415649","3.63851833343506","3.63851833343506","3.63851833343506","3.63851833343506","3.63851833343506","3.81363439559937","3.81363439559937","3.75526213645935","3.63851833343506","3.63851833343506","3.63851833343506","3.63851833343506","3.69689035415649","3.69689035415649","3.69689035415649","3.69689035415649","3.69689035415649","3.81363439559937","3.81363439559937","3.81363439559937","3.85254883766174","3.87200617790222","4.00820732116699","4.08603668212891","4.08603668212891","4.08603668212891","4.08603668212891","4.08603668212891","4.08603668212891","4.08603668212891","3.98874998092651","3.98874998092651","3.98874998092651","3.98874998092651","3.75526213645935","3.59960389137268","3.50231719017029","3.50231719017029","3.50231719017029","3.46340250968933","3.46340250968933","3.59960389137268","3.59960389137268","3.69689035415649","3.69689035415649","3.69689035415649","3.69689035415649","3.69689035415649","3.44394516944885","3.59960389137268","3.59960389137268","3.59960389137268","3.4050304889679","3.52177429199219","3.52177429199219","3.52177429199219","3.52177429199219","3.385573387146","3.385573387146","3.385573387146","2.91859769821167","2.87968301773071","2.87968301773071","2.84076833724976","2.91859769821167","2.91859769821167","2.82131123542786","2.87968301773071","2.87968301773071","2.87968301773071","2.87968301773071","2.80185389518738","2.80185389518738","2.78239631652832","2.76293897628784","2.70456719398499","2.70456719398499","2.68510985374451","2.66565251350403","2.66565251350403","2.68510985374451","2.68510985374451","2.70456719398499","2.70456719398499","2.70456719398499","2.70456719398499","2.70456719398499","2.70456719398499","2.91859769821167","2.72402453422546","2.70456719398499","2.66565251350403","2.66565251350403","2.76293897628784","2.76293897628784","2.58782339096069","2.60728073120117","2.64619517326355","2.66565251350403","2.41270732879639","2.43216466903687","2.41270732879639","2.41270732879639","2.43216466903687","2.43216466903687","2.37379288673401","2.29596328735352","2.43216466903687","2.33487820625305","2.31542086601257","2.23759150505066","2.23759150505066","2.12084770202637","2.21813416481018","2.04301834106445","2.02356100082397","2.02356100082397","2.02356100082397","1.98464632034302","2.25704884529114","2.14030504226685","2.06247544288635","2.06247544288635","2.06247544288635","2.27650618553162","2.25704884529114","2.41270732879639","2.43216466903687","2.43216466903687","2.43216466903687","2.29596328735352","2.21813416481018","2.14030504226685","2.15976238250732","2.31542086601257","2.31542086601257","2.17921948432922","2.15976238250732","2.15976238250732","2.23759150505066","2.27650618553162","2.27650618553162","2.14030504226685","2.17921948432922","2.17921948432922","2.21813416481018","2.15976238250732","2.14030504226685","2.15976238250732","2.21813416481018","2.1986768245697","2.31542086601257","2.52945113182068","3.77471971511841","3.30774402618408","3.0937135219574","2.996426820755","2.91859769821167","2.76293897628784","2.82131123542786","2.89914035797119","3.01588416099548","2.95751214027405","2.87968301773071","2.91859769821167","2.95751214027405","2.95751214027405","2.86022567749023","2.87968301773071","2.91859769821167","2.91859769821167","2.95751214027405","2.91859769821167","2.87968301773071","2.97696948051453","3.01588416099548","3.13262820243835","3.13262820243835","3.21045732498169","3.19099998474121")

# 1) Append the new raw OHLCV row (482) right after the last existing row (481).
#    Clone A481's date style (custom "yyyy-mm-dd hh:mm:ss" format) onto A482.
$ws.Cells.Item(481, 1).Copy()
$ws.Cells.Item(482, 1).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(482, 1).Value = 45450.2916666667   # date
$ws.Cells.Item(482, 2).Value = 44700              # volume
$ws.Cells.Item(482, 3).Value = 3.38000011444092   # high
$ws.Cells.Item(482, 4).Value = 3.25999999046326   # low
$ws.Cells.Item(482, 5).Value = 3.29999995231628   # open
$ws.Cells.Item(482, 6).Value = 3.27999997138977   # close
$ws.Cells.Item(482, 8).Value = "IMP.MI"           # ticker

# 2) Rewrite column G (adj_close) for every data row (2..482). These values
#    are decimal strings stored as text in the sheet (matching the source
#    file's convention), so a throwaway cell is pre-formatted as Text ("@")
#    and its value is copy/pasted as values-only into each target cell -
#    this keeps the numeric-looking text from being auto-converted back to
#    a number, without leaving any custom formatting behind on the data cells.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.NumberFormat = "@"

for ($i = 0; $i -lt $adjClose.Length; $i++) {
    $row = $i + 2
    $scratch.Value = $adjClose[$i]
    $scratch.Copy()
    $ws.Cells.Item($row, 7).PasteSpecial(-4163)   # xlPasteValues
}

$scratch.Clear()

Write-Output "Updated $($adjClose.Length) adj_close cells and appended row 482."
